$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain identical data tables and both
# need the same updates applied (mirrors the diff, which touches the
# sheet1.xml and sheet4.xml parts identically).
$sheetNames = @("展览", "全部类型")

# Map of row number -> updates for columns F (想去人数) and/or G (最低票价)
$updates = @{
    3  = @{ F = 3124; G = 75 }
    7  = @{ F = 1731 }
    11 = @{ F = 7 }
    12 = @{ F = 1419 }
    15 = @{ F = 354 }
    16 = @{ F = 72 }
    17 = @{ F = 11 }
    19 = @{ F = 63 }
    23 = @{ F = 114 }
    24 = @{ F = 3338 }
    25 = @{ F = 402 }
    26 = @{ F = 200 }
    27 = @{ F = 449 }
    28 = @{ F = 27 }
    29 = @{ F = 17 }
    31 = @{ F = 1047 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $cellRef = "$col$row"
            $ws.Range($cellRef).Value = $cols[$col]
        }
    }
}
